# Generate Report for Handoff
# Updates the localization-status report: the four files that were queued
# with "low" priority have now been handed off, so their Priority flips to
# "ht" and their Latest Handoff Datetime is refreshed on both the zh-cn and
# de-de sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-21 12:38:53"
}

$dede = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-21 12:38:57"
}

# The de-de "Latest Handoff Datetime" for these same rows previously shared
# its string value with the Overview sheet's "Latest HO Xliff Generate Date"
# column, so refresh that column in lockstep to keep both in sync.
$overview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-08-21 12:38:57"
}
